# Applies updated mass-flow values (new input files) to both worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Output_flows")
$ws.Range("C2").Value = [double]"1.133510123900368E-29"
$ws.Range("E2").Value = [double]"3.313267739150016E-29"
$ws.Range("H2").Value = [double]"3.306203202392766E-44"
$ws.Range("C3").Value = [double]"1.714860579195604E-30"
$ws.Range("D3").Value = [double]"1.694534831110598E-31"
$ws.Range("E3").Value = [double]"5.012564170700179E-30"
$ws.Range("H3").Value = [double]"5.001876400613398E-45"
$ws.Range("C4").Value = [double]"7.464415303555643E-31"
$ws.Range("D4").Value = [double]"7.375941740805077E-31"
$ws.Range("E4").Value = [double]"2.181860214162712E-30"
$ws.Range("H4").Value = [double]"2.177208060187951E-45"
$ws.Range("C5").Value = [double]"1.238912770834886E-31"
$ws.Range("D5").Value = [double]"1.224228294916093E-30"
$ws.Range("E5").Value = [double]"3.621361317094852E-31"
$ws.Range("H5").Value = [double]"3.613639864393157E-46"
$ws.Range("C6").Value = [double]"3.413061473587765E-31"
$ws.Range("D6").Value = [double]"3.372607439859073E-29"
$ws.Range("E6").Value = [double]"9.976431823354519E-31"
$ws.Range("H6").Value = [double]"9.95516011370984E-46"
$ws.Range("C7").Value = [double]"2.060554648301047E-22"
$ws.Range("E7").Value = [double]"6.023033316614253E-21"
$ws.Range("H7").Value = [double]"6.010191028092687E-36"
$ws.Range("C8").Value = [double]"4.532909392242179E-23"
$ws.Range("D8").Value = [double]"4.480674687260957E-25"
$ws.Range("E8").Value = [double]"1.324976472387133E-21"
$ws.Range("H8").Value = [double]"1.322151362637907E-36"
$ws.Range("C9").Value = [double]"3.012955742384946E-23"
$ws.Range("D9").Value = [double]"2.97724503896003E-24"
$ws.Range("E9").Value = [double]"8.806916542024883E-22"
$ws.Range("H9").Value = [double]"8.788138468374511E-37"
$ws.Range("C10").Value = [double]"1.458839299422834E-23"
$ws.Range("D10").Value = [double]"1.441548098241648E-23"
$ws.Range("E10").Value = [double]"4.264209984071335E-22"
$ws.Range("H10").Value = [double]"4.255117851909135E-37"
$ws.Range("C11").Value = [double]"8.732879137029968E-22"
$ws.Range("D11").Value = [double]"8.629370838137006E-21"
$ws.Range("E11").Value = [double]"2.552634167488112E-20"
$ws.Range("H11").Value = [double]"2.547191450713084E-35"
$ws.Range("C12").Value = [double]"2.805123917720582E-28"
$ws.Range("E12").Value = [double]"4.09971044145638E-27"
$ws.Range("H12").Value = [double]"4.090969054587602E-42"
$ws.Range("C13").Value = [double]"6.652664023666119E-29"
$ws.Range("D13").Value = [double]"4.93035893611825E-29"
$ws.Range("E13").Value = [double]"9.722920256402719E-28"
$ws.Range("H13").Value = [double]"9.702189083148129E-43"
$ws.Range("C14").Value = [double]"4.237491458612639E-29"
$ws.Range("D14").Value = [double]"1.256179708181464E-28"
$ws.Range("E14").Value = [double]"6.193126752337267E-28"
$ws.Range("H14").Value = [double]"6.179921791244863E-43"
$ws.Range("C15").Value = [double]"6.190931699776104E-29"
$ws.Range("D15").Value = [double]"1.559975854620591E-27"
$ws.Range("E15").Value = [double]"9.048094870810476E-28"
$ws.Range("H15").Value = [double]"9.028802557653097E-43"
$ws.Range("C16").Value = [double]"1.920184818978827E-29"
$ws.Range("D16").Value = [double]"4.753050681785765E-27"
$ws.Range("E16").Value = [double]"2.806365060082766E-28"
$ws.Range("H16").Value = [double]"2.800381339272355E-43"
$ws.Range("C17").Value = [double]"8.887315919527002E-21"
$ws.Range("E17").Value = [double]"5.19555255889208E-19"
$ws.Range("H17").Value = [double]"5.184474621666297E-34"
$ws.Range("C18").Value = [double]"1.950626549100046E-21"
$ws.Range("D18").Value = [double]"5.782590393063401E-23"
$ws.Range("E18").Value = [double]"1.140342354248047E-19"
$ws.Range("H18").Value = [double]"1.1379109206569E-34"
$ws.Range("C19").Value = [double]"1.340817957344596E-21"
$ws.Range("D19").Value = [double]"1.589911061048499E-22"
$ws.Range("E19").Value = [double]"7.838463527535912E-20"
$ws.Range("H19").Value = [double]"7.821750385686154E-35"
$ws.Range("C20").Value = [double]"3.84528841497945E-20"
$ws.Range("D20").Value = [double]"3.8757055475326E-20"
$ws.Range("E20").Value = [double]"2.247967580428696E-18"
$ws.Range("H20").Value = [double]"2.2431744725813E-33"
$ws.Range("C21").Value = [double]"7.28507801551406E-20"
$ws.Range("D21").Value = [double]"7.213127546113547E-19"
$ws.Range("E21").Value = [double]"4.258879291335796E-18"
$ws.Range("H21").Value = [double]"4.249798525256189E-33"

$ws = $wb.Worksheets.Item("Input_flows")
$ws.Range("C2").Value = [double]"1.95867390043894E-29"
$ws.Range("C3").Value = [double]"1.224324327194767E-30"
$ws.Range("C4").Value = [double]"4.844281975649409E-32"
$ws.Range("C5").Value = [double]"2.395198377953042E-32"
$ws.Range("C6").Value = [double]"3.506502372828496E-29"
$ws.Range("C7").Value = [double]"1.779432161708033E-22"
$ws.Range("C8").Value = [double]"7.197628934176793E-23"
$ws.Range("C9").Value = [double]"4.797827665511482E-23"
$ws.Range("C10").Value = [double]"2.395633047692827E-23"
$ws.Range("C11").Value = [double]"3.502900042672115E-20"
$ws.Range("C12").Value = [double]"4.989268435796452E-30"
$ws.Range("C13").Value = [double]"3.604966370603483E-31"
$ws.Range("C14").Value = [double]"5.321535305480493E-33"
$ws.Range("C15").Value = [double]"2.289042124610113E-27"
$ws.Range("C16").Value = [double]"5.052889035983831E-27"
$ws.Range("C17").Value = [double]"1.32088838608534E-22"
$ws.Range("C18").Value = [double]"1.486803801464846E-23"
$ws.Range("C19").Value = [double]"1.757782607967625E-24"
$ws.Range("C20").Value = [double]"2.289111882323251E-18"
$ws.Range("C21").Value = [double]"5.053042826102296E-18"
